# Swap the report rows for the two localization files
# "7ede8c4c-a33d-4d09-bfc7-0c757559c58f" and "3e973b86-4134-4cf4-a4f6-4160150e0136"
# (row 4 <-> row 5) across the Overview, zh-cn and de-de sheets, including the
# hyperlink display text (the underlying hyperlink target / rId stays attached
# to the row position, only the visible text moves).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A (File Name), B (Path And Name, hyperlinked),
# E/F (zh-cn / de-de status) and G (Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "3e973b86-4134-4cf4-a4f6-4160150e0136.md"
$wsOverview.Range("B4").Value = "e2e\3e973b86-4134-4cf4-a4f6-4160150e0136.md"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"
$wsOverview.Range("G4").Value = "2016-09-06 08:55:42"

$wsOverview.Range("A5").Value = "7ede8c4c-a33d-4d09-bfc7-0c757559c58f.md"
$wsOverview.Range("B5").Value = "e2e\7ede8c4c-a33d-4d09-bfc7-0c757559c58f.md"
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"
$wsOverview.Range("G5").Value = "2016-09-06 08:54:49"

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$4') {
        $h.TextToDisplay = "e2e\3e973b86-4134-4cf4-a4f6-4160150e0136.md"
    }
    elseif ($addr -eq '$B$5') {
        $h.TextToDisplay = "e2e\7ede8c4c-a33d-4d09-bfc7-0c757559c58f.md"
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn": columns A (Source File Name), C (Status), G (Latest Handoff
# File) and H (Latest Handoff Datetime)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "3e973b86-4134-4cf4-a4f6-4160150e0136.md"
$wsZhCn.Range("C4").Value = "In Translation"
$wsZhCn.Range("G4").Value = "3e973b86-4134-4cf4-a4f6-4160150e0136.38a31fdee6349160df87bce430e057f82dbcfed2.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-09-06 08:55:37"

$wsZhCn.Range("A5").Value = "7ede8c4c-a33d-4d09-bfc7-0c757559c58f.md"
$wsZhCn.Range("C5").Value = "In Translation"
$wsZhCn.Range("G5").Value = "7ede8c4c-a33d-4d09-bfc7-0c757559c58f.7a928b72346d3c7982f8a6d33ff6b689039ac6dc.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-09-06 08:54:43"

foreach ($h in $wsZhCn.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$4') {
        $h.TextToDisplay = "3e973b86-4134-4cf4-a4f6-4160150e0136.md"
    }
    elseif ($addr -eq '$A$5') {
        $h.TextToDisplay = "7ede8c4c-a33d-4d09-bfc7-0c757559c58f.md"
    }
}

# ---------------------------------------------------------------------------
# Sheet "de-de": columns A (Source File Name), C (Status), G (Latest Handoff
# File) and H (Latest Handoff Datetime)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "3e973b86-4134-4cf4-a4f6-4160150e0136.md"
$wsDeDe.Range("C4").Value = "In Translation"
$wsDeDe.Range("G4").Value = "3e973b86-4134-4cf4-a4f6-4160150e0136.38a31fdee6349160df87bce430e057f82dbcfed2.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-09-06 08:55:42"

$wsDeDe.Range("A5").Value = "7ede8c4c-a33d-4d09-bfc7-0c757559c58f.md"
$wsDeDe.Range("C5").Value = "In Translation"
$wsDeDe.Range("G5").Value = "7ede8c4c-a33d-4d09-bfc7-0c757559c58f.7a928b72346d3c7982f8a6d33ff6b689039ac6dc.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-09-06 08:54:49"

foreach ($h in $wsDeDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$4') {
        $h.TextToDisplay = "3e973b86-4134-4cf4-a4f6-4160150e0136.md"
    }
    elseif ($addr -eq '$A$5') {
        $h.TextToDisplay = "7ede8c4c-a33d-4d09-bfc7-0c757559c58f.md"
    }
}
